$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 138 and 152: the name entries are removed (cells cleared) first, so the
# two retired shared strings ("سامي احلو" / "الزهارنة ألمنيوم") drop out of the
# shared-string table before any new strings are appended. Clear() (rather than
# ClearContents()) drops the <c> element entirely, matching a truly blank cell.
$ws.Range("A138").Clear()
$ws.Range("A152").Clear()

# --- Column A text updates ---
# New shared strings must be introduced in this order (A32, then A110, then
# G1) so they land at the expected shared-string-table indices.
# Row 32: string text actually changes
$ws.Range("A32").Value = "صندوق2 صندوق احمد صندوق2 احمد"

# Row 110: previously blank, now gets a new name
$ws.Range("A110").Value = "شرين مشتهى"

# --- New header column G (يورو) ---
$ws.Range("G1").Value = "يورو"

# --- Numeric value updates (columns B/C/D) ---
$ws.Range("B11").Value = -73445
$ws.Range("C11").Value = -26721

$ws.Range("C16").Value = -1805

$ws.Range("D32").Value = -27712

$ws.Range("B33").Value = 32000

$ws.Range("B36").Value = 4800

$ws.Range("C53").Value = 6500

$ws.Range("D58").Value = 6762

$ws.Range("C80").Value = -13079

$ws.Range("C94").Value = -2820

$ws.Range("B104").Value = -4400

$ws.Range("C110").Value = -600

$ws.Range("C116").Value = 1000

$ws.Range("B117").Value = 132000
$ws.Range("C117").Value = -100

$ws.Range("C120").Value = 7300

$ws.Range("C138").Value = 0

$ws.Range("C141").Value = -139830

$ws.Range("C143").Value = -1080

$ws.Range("C147").Value = -53384

$ws.Range("C152").Value = 0

# --- Sheet view: scroll / selection state ---
$ws.Range("F156").Select()
